# Complete test for table_CalendarReturns macro
#
# The old row describing the (incomplete) table_CalendarReturns test is
# removed from its original position; the three rows that followed it
# shift up by one; and the CalendarReturns test is completed and moved
# to the bottom of the table as two distinct rows (simple vs. compound
# returns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 73-75 shift up to 72-74 (their content is unchanged, only their
# row position moves, because the old CalendarReturns row that used to
# sit at row 72 is gone).
$ws.Range("A72").Value = "Table_InformationRatio1"
$ws.Range("B72").Value = "Test Information Ratio Table with scale=252"
$ws.Range("C72").Value = "table_InformationRatio_test1"

$ws.Range("A73").Value = "Table_Stats"
$ws.Range("B73").Value = "Test stats table"
$ws.Range("C73").Value = "table_stats_test"

$ws.Range("A74").Value = "Table_Correlation"
$ws.Range("B74").Value = "Test correlation table"
$ws.Range("C74").Value = "table_correlation_test"

# The completed CalendarReturns tests land at the end of the table as
# two rows: one for simple returns, one for compound returns.
$ws.Range("A75").Value = "Table_CalendarReturns1"
$ws.Range("B75").Value = "Test Calendar Returns for simple returns with digits=6"
$ws.Range("C75").Value = "table_CalendarReturns_test1"

$ws.Range("A76").Value = "Table_CalendarReturns2"
$ws.Range("B76").Value = "Test Calendar Returns for compound returns with digits=8"
$ws.Range("C76").Value = "table_CalendarReturns_test2"

# Keep the view roughly where the author left it: scrolled so the new
# rows are visible, with the cell just past the last data row selected.
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("C77").Select()
